# Re-apply the bold / not-italic / not-strikethrough character formatting on
# the four heading-style runs ("Starting with M2Doc", "Template user",
# "Template developper", "Integrator"). The formatting values themselves are
# unchanged (bold on, italic off, strike off) - this simply forces the
# run properties to be re-serialized by the current engine (equivalent to
# the POI 4.1.0 -> 5.2.3 upgrade that changed how w:b/w:i/w:strike on-off
# values are written), without touching the font size (w:sz) that must stay
# as-is.

$d = $word.ActiveDocument

$targets = @("Starting with M2Doc", "Template user", "Template developper", "Integrator")

foreach ($t in $targets) {
    $rng = $d.Content
    $found = $rng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = $true
        $rng.Font.Italic = $false
        $rng.Font.StrikeThrough = $false
    }
}
